$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43: Prompt_id 1 - Results @ 1, 0-shot
$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 101
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 101
$ws.Range("G43").Value = 1

# Row 48: Prompt_id 1 - Results @ 2, 0-shot
$ws.Range("B48").Value = 0
$ws.Range("C48").Value = 2
$ws.Range("D48").Value = 99
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 200
$ws.Range("G48").Value = 0.9901

# Row 53: Prompt_id 1 - Results @ 5, 0-shot
$ws.Range("B53").Value = 1
$ws.Range("C53").Value = 3
$ws.Range("D53").Value = 97
$ws.Range("E53").Value = 0.0099
$ws.Range("F53").Value = 499
$ws.Range("G53").Value = 0.98812

# Update the view selection to match the saved workbook state
$ws.Range("B76").Select()
